# Update the two "_694" model-run labels to reflect the Plan Amendment
# (2035_06_694 -> 2035_06_694_Amd1, 2040_06_694 -> 2040_06_694_Amd1).
# Write the "_694" row first so the new shared-string entries land in the
# order 2040_06_694_Amd1, 2035_06_694_Amd1 (matching the authored workbook).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(18, 2).Value2 = "2040_06_694_Amd1"
$ws.Cells.Item(13, 2).Value2 = "2035_06_694_Amd1"

# Column B needs to be widened to fit the longer "_Amd1" labels.
$ws.Columns.Item(2).ColumnWidth = 62/3

# Move/restore the active selection to B22.
$ws.Range("B22").Select() | Out-Null
